$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume-change (E) columns for rows with new market data
$ws.Range("D2").Value = "'29.117.89"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").Value = "'1.926.49"
$ws.Range("E3").Value = "  +1.78%  "
$ws.Range("D4").Value = "'0.9975"
$ws.Range("E4").Value = "  -0.80%  "
$ws.Range("D5").Value = "'325.14"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("D6").Value = "'0.9975"
$ws.Range("E6").Value = "  -0.79%  "
$ws.Range("D7").Value = "'0.4616"
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").Value = "'0.3882"
$ws.Range("E8").Value = "  -0.83%  "
$ws.Range("D9").Value = "'0.07847"
$ws.Range("E9").Value = "  -0.32%  "
$ws.Range("D10").Value = "'0.9941"
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("D11").Value = "'22.04"
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("D12").Value = "'1.913.37"
$ws.Range("E12").Value = "  +1.82%  "
$ws.Range("D13").Value = "'5.787"
$ws.Range("E13").Value = "  +1.35%  "
$ws.Range("D14").Value = "'7.052"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").Value = "'0.07046"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").Value = "'87.51"
$ws.Range("E16").Value = "  -0.79%  "
$ws.Range("D17").Value = "'0.9988"
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("D18").Value = "'0.000009926"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("D19").Value = "'17.08"
$ws.Range("E19").Value = "  +0.76%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "'29.115.92"
$ws.Range("E21").Value = "  +0.85%  "
$ws.Range("D22").Value = "'5.402"
$ws.Range("E22").Value = "  +1.68%  "
$ws.Range("E23").Value = "  +1.46%  "
$ws.Range("D24").Value = "'2.132.14"
$ws.Range("E24").Value = "  +1.11%  "
$ws.Range("D25").Value = "'2.087"
$ws.Range("E25").Value = "  +1.09%  "
$ws.Range("D26").Value = "'155.90"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").Value = "'19.42"
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("D28").Value = "'5.871"
$ws.Range("E28").Value = "  -1.41%  "
$ws.Range("D29").Value = "'118.40"
$ws.Range("E29").Value = "  +0.65%  "
$ws.Range("D30").Value = "'1.884"
$ws.Range("E30").Value = "  -2.46%  "
$ws.Range("D31").Value = "'0.09320"
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("D32").Value = "'0.8855"
$ws.Range("E32").Value = "  -3.09%  "
$ws.Range("D33").Value = "'5.207"
$ws.Range("E33").Value = "  -1.92%  "
$ws.Range("D34").Value = "'1.320"
$ws.Range("E34").Value = "  -1.19%  "
$ws.Range("D35").Value = "'3.117"
$ws.Range("E35").Value = "  -5.37%  "
$ws.Range("D36").Value = "'0.05764"
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("D37").Value = "'1.171"
$ws.Range("E37").Value = "  -1.39%  "
$ws.Range("D38").Value = "'0.02092"
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("D39").Value = "'0.9966"
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("D40").Value = "'7.664"
$ws.Range("E40").Value = "  -1.37%  "
$ws.Range("D41").Value = "'0.5690"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Value = "'0.1811"
$ws.Range("E42").Value = "  +1.74%  "
$ws.Range("D43").Value = "'0.000003052"
$ws.Range("E43").Value = "  +97.06%  "
$ws.Range("D44").Value = "'9.725"
$ws.Range("E44").Value = "  -0.69%  "
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").Value = "'2.217"
$ws.Range("E46").Value = "  -1.77%  "
$ws.Range("D47").Value = "'0.5326"
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("D48").Value = "'0.06927"
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("D51").Value = "'112.82"
$ws.Range("E51").Value = "  +0.24%  "

# Rows 49 and 50 swapped rank order (NEARProtocol moved above MXToken)
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'1.843"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").Value = "'2.574"
$ws.Range("E50").Value = "  +2.13%  "
